$wb = $excel.ActiveWorkbook

# Insert the new "World" worksheet after the last existing sheet (Tabelle3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "World"

$data = @(
    @('Year', 'Value', 'Location'),
    @('"2000"', '"61.25434"', '"ALLAREA"'),
    @('"2000"', '"86.17919"', '"URBAN"'),
    @('"2000"', '"39.4361"', '"RURAL"'),
    @('"2001"', '"61.70454"', '"ALLAREA"'),
    @('"2001"', '"86.13555"', '"URBAN"'),
    @('"2001"', '"39.92612"', '"RURAL"'),
    @('"2002"', '"86.34688"', '"URBAN"'),
    @('"2002"', '"62.83692"', '"ALLAREA"'),
    @('"2002"', '"41.45364"', '"RURAL"'),
    @('"2003"', '"86.29749"', '"URBAN"'),
    @('"2003"', '"42.55228"', '"RURAL"'),
    @('"2003"', '"63.60941"', '"ALLAREA"'),
    @('"2004"', '"63.93015"', '"ALLAREA"'),
    @('"2004"', '"42.79121"', '"RURAL"'),
    @('"2004"', '"86.24792"', '"URBAN"'),
    @('"2005"', '"43.03233"', '"RURAL"'),
    @('"2005"', '"64.24246"', '"ALLAREA"'),
    @('"2005"', '"86.18039"', '"URBAN"'),
    @('"2006"', '"86.12835"', '"URBAN"'),
    @('"2006"', '"43.29615"', '"RURAL"'),
    @('"2006"', '"64.56688"', '"ALLAREA"'),
    @('"2007"', '"86.07311"', '"URBAN"'),
    @('"2007"', '"64.88221"', '"ALLAREA"'),
    @('"2007"', '"43.55946"', '"RURAL"'),
    @('"2008"', '"43.8296"', '"RURAL"'),
    @('"2008"', '"65.19863"', '"ALLAREA"'),
    @('"2008"', '"86.01056"', '"URBAN"'),
    @('"2009"', '"44.8021"', '"RURAL"'),
    @('"2009"', '"65.85203"', '"ALLAREA"'),
    @('"2009"', '"85.94473"', '"URBAN"'),
    @('"2010"', '"85.87328"', '"URBAN"'),
    @('"2010"', '"66.49626"', '"ALLAREA"'),
    @('"2010"', '"45.78517"', '"RURAL"'),
    @('"2011"', '"46.81956"', '"RURAL"'),
    @('"2011"', '"85.77981"', '"URBAN"'),
    @('"2011"', '"67.12907"', '"ALLAREA"'),
    @('"2012"', '"47.85402"', '"RURAL"'),
    @('"2012"', '"67.74305"', '"ALLAREA"'),
    @('"2012"', '"85.68166"', '"URBAN"'),
    @('"2013"', '"48.89293"', '"RURAL"'),
    @('"2013"', '"68.34628"', '"ALLAREA"'),
    @('"2013"', '"85.57937"', '"URBAN"'),
    @('"2014"', '"49.93764"', '"RURAL"'),
    @('"2014"', '"85.4778"', '"URBAN"'),
    @('"2014"', '"68.94266"', '"ALLAREA"'),
    @('"2015"', '"50.98648"', '"RURAL"'),
    @('"2015"', '"85.37737"', '"URBAN"'),
    @('"2015"', '"69.53274"', '"ALLAREA"'),
    @('"2016"', '"85.26961"', '"URBAN"'),
    @('"2016"', '"52.0242"', '"RURAL"'),
    @('"2016"', '"70.10356"', '"ALLAREA"'),
    @('"2017"', '"53.03716"', '"RURAL"'),
    @('"2017"', '"70.64488"', '"ALLAREA"'),
    @('"2017"', '"85.14748"', '"URBAN"')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    for ($j = 0; $j -lt 3; $j++) {
        $newSheet.Cells.Item($i + 1, $j + 1).Value = $data[$i][$j]
    }
}

# Match the captured selection/scroll state on the new sheet.
$newSheet.Range("I36").Select()
